# General Fixes: update the "Expected status" comments in column D to also
# record the actual observed status ("Ожидаемый статус - X/Фактический статус - Y"),
# and move the active selection to D47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = "Ожидаемый статус - 200/Фактический статус - 200"
$ws.Range("D6").Value  = "Ожидаемый статус - 200/Фактический статус - 200"
$ws.Range("D7").Value  = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D8").Value  = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D9").Value  = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D10").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D11").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D12").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D13").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D14").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D15").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D16").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D17").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D18").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D19").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D20").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D21").Value = "Ожидаемый статус - 200/Фактический статус - 200"
$ws.Range("D22").Value = "Ожидаемый статус - 200/Фактический статус - 200"
$ws.Range("D23").Value = "Ожидаемый статус - 200/Фактический статус - 200"
$ws.Range("D24").Value = "Ожидаемый статус - 200/Фактический статус - 200"
$ws.Range("D25").Value = "Ожидаемый статус - 200/Фактический статус - 200"
$ws.Range("D26").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D27").Value = "Ожидаемый статус - 400/Фактический статус - 400"
$ws.Range("D28").Value = "Ожидаемый статус - 200/Фактический статус - 400"
$ws.Range("D29").Value = "Ожидаемый статус - 200/Фактический статус - 400"
$ws.Range("D30").Value = "Ожидаемый статус - 400/Фактический статус - 400"
$ws.Range("D31").Value = "Ожидаемый статус - 200/Фактический статус - 400"
$ws.Range("D32").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D33").Value = "Ожидаемый статус - 400/Фактический статус - 400"
$ws.Range("D34").Value = "Ожидаемый статус - 400/Фактический статус - 400"
$ws.Range("D35").Value = "Ожидаемый статус - 400/Фактический статус - 400"
$ws.Range("D36").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D37").Value = "Ожидаемый статус - 400/Фактический статус - 200"
$ws.Range("D39").Value = "Ожидаемый статус - 200/Фактический статус - 200"
$ws.Range("D41").Value = "Ожидаемый статус - 401/Фактический статус - 401"
$ws.Range("D42").Value = "Ожидаемый статус - 401/Фактический статус - 401"
$ws.Range("D43").Value = "Ожидаемый статус - 401/Фактический статус - 401"
$ws.Range("D44").Value = "Ожидаемый статус - 401/Фактический статус - 401"

$ws.Range("D47").Select()
